$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.179.29"
$ws.Range("E2").Value = "  -1.16%  "

$ws.Range("D3").Value = "1.659.06"
$ws.Range("E3").Value = "  -1.11%  "

$ws.Range("E4").Value = "  +0.42%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.48%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5199"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.41%  "

$ws.Range("E7").Value = "  +0.40%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2628"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.58%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06240"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.68%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.69"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.68%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07709"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.08%  "

$ws.Range("D12").Value = "1.657.09"
$ws.Range("E12").Value = "  -1.22%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.395"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.68%  "

$ws.Range("D14").Value = "1.884.29"
$ws.Range("E14").Value = "  -1.18%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5416"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.21%  "

$ws.Range("D16").Value = "0.0₅8074"
$ws.Range("E16").Value = "  -2.97%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.26"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.23%  "

$ws.Range("D18").Value = "26.195.37"
$ws.Range("E18").Value = "  -1.21%  "

$ws.Range("E19").Value = "  +0.47%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.606"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.80%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "190.88"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.43%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.01"
$ws.Range("D22").Style = "Normal"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.043"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.70%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.007"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.53%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "139.81"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.82%  "

$ws.Range("E26").Value = "  -4.92%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.124"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.01%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.98"
$ws.Range("D28").Style = "Normal"

$ws.Range("E29").Value = "  -2.66%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05970"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.00%  "

$ws.Range("E31").Value = "  -0.20%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.561"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.61%  "

$ws.Range("E33").Value = "  -6.47%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.612"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.14%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9599"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.98%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.419"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.08%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.777"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.47%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5642"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -8.18%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.006"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.53%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01587"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.87%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8558"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.95%  "

$ws.Range("D43").Value = "1.018.19"
$ws.Range("E43").Value = "  -7.12%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.93"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.65%  "

$ws.Range("D45").Value = "1.799.30"
$ws.Range("E45").Value = "  -1.29%  "

$ws.Range("E46").Value = "  +3.59%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "56.55"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.43%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.012"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.95%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.916"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.89%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05176"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.65%  "

$ws.Range("E51").Value = "  -0.56%  "

